$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns (D:E) for the two newest quarters; existing data shifts right
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats / styles from column F (the old column D, now shifted) into D:E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 36/37 and 78/79 are separator/section-title rows with no data columns;
# undo the formatting paste there so no stray D:E cells are introduced
$ws.Range("D36:E37").Clear()
$ws.Range("D78:E79").Clear()

# Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30)
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 8900
$ws.Cells.Item(8,5).Value = 7000
$ws.Cells.Item(9,4).Value = 1900
$ws.Cells.Item(9,5).Value = 2100
$ws.Cells.Item(10,4).Value = 7000
$ws.Cells.Item(10,5).Value = 4900
$ws.Cells.Item(12,4).Value = "NA"
$ws.Cells.Item(12,5).Value = "NA"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(15,4).Value = 1000
$ws.Cells.Item(15,5).Value = 1000
$ws.Cells.Item(17,4).Value = 12700
$ws.Cells.Item(17,5).Value = 12300
$ws.Cells.Item(18,4).Value = -3800
$ws.Cells.Item(18,5).Value = -5300
$ws.Cells.Item(20,4).Value = -1100
$ws.Cells.Item(20,5).Value = -100
$ws.Cells.Item(21,4).Value = -4000
$ws.Cells.Item(21,5).Value = -4400
$ws.Cells.Item(22,4).Value = 1600
$ws.Cells.Item(22,5).Value = 1700
$ws.Cells.Item(23,4).Value = -6500
$ws.Cells.Item(23,5).Value = -7100
$ws.Cells.Item(24,4).Value = -100
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = -6400
$ws.Cells.Item(26,5).Value = -7100
$ws.Cells.Item(27,4).Value = -6400
$ws.Cells.Item(27,5).Value = -7100
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = 1100
$ws.Cells.Item(32,5).Value = 100
$ws.Cells.Item(33,4).Value = -6400
$ws.Cells.Item(33,5).Value = -7100
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = -6400
$ws.Cells.Item(35,5).Value = -7100
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 15600
$ws.Cells.Item(41,5).Value = 16800
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,4).Value = 7700
$ws.Cells.Item(43,5).Value = 8600
$ws.Cells.Item(44,4).Value = 4200
$ws.Cells.Item(44,5).Value = 4900
$ws.Cells.Item(45,4).Value = 2800
$ws.Cells.Item(45,5).Value = 3200
$ws.Cells.Item(46,4).Value = 30300
$ws.Cells.Item(46,5).Value = 33500
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(48,4).Value = 500
$ws.Cells.Item(48,5).Value = 500
$ws.Cells.Item(49,4).Value = 26800
$ws.Cells.Item(49,5).Value = 27700
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 2000
$ws.Cells.Item(52,5).Value = 300
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 59600
$ws.Cells.Item(54,5).Value = 62100
$ws.Cells.Item(57,4).Value = 4200
$ws.Cells.Item(57,5).Value = 3400
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(59,4).Value = 5200
$ws.Cells.Item(59,5).Value = 4900
$ws.Cells.Item(60,4).Value = 9400
$ws.Cells.Item(60,5).Value = 8300
$ws.Cells.Item(61,4).Value = 41500
$ws.Cells.Item(61,5).Value = 40700
$ws.Cells.Item(62,4).Value = 1800
$ws.Cells.Item(62,5).Value = 2900
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 52700
$ws.Cells.Item(66,5).Value = 51900
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = -409700
$ws.Cells.Item(72,5).Value = -403300
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 6900
$ws.Cells.Item(76,5).Value = 10100
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = -6400
$ws.Cells.Item(81,5).Value = -7100
$ws.Cells.Item(83,4).Value = 1000
$ws.Cells.Item(83,5).Value = 1000
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = -4400
$ws.Cells.Item(89,5).Value = -7500
$ws.Cells.Item(91,4).Value = 0
$ws.Cells.Item(91,5).Value = 0
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = 200
$ws.Cells.Item(94,5).Value = 200
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = 3000
$ws.Cells.Item(100,5).Value = 2200
$ws.Cells.Item(101,4).Value = -100
$ws.Cells.Item(101,5).Value = 0
$ws.Cells.Item(102,4).Value = -1300
$ws.Cells.Item(102,5).Value = -5100

# Restated historical figures that changed as part of this data refresh
$ws.Cells.Item(14,6).Value = -18500
$ws.Cells.Item(17,6).Value = -2700
$ws.Cells.Item(18,6).Value = 8900
$ws.Cells.Item(20,6).Value = -1700
$ws.Cells.Item(32,6).Value = 1700
$ws.Cells.Item(91,6).Value = -100
$ws.Cells.Item(91,9).Value = 0
$ws.Cells.Item(91,10).Value = "NA"
$ws.Cells.Item(101,8).Value = 100
$ws.Cells.Item(101,9).Value = 300
